# Supporting-documentation template ("Tribals") update:
#   Reference #_reference_num_, MP(s)_mps_
# becomes
#   Reference #_reference_num_, MP(s) - _mps_
# with the (hidden) "_GoBack" bookmark moving from just after "County" in
# the Location paragraph to sit right before "_mps_" in the Reference
# paragraph.

$d = $word.ActiveDocument

# 1) Insert " - " right before the "_mps_" placeholder on the
#    "Reference #..." line. Doing this as a Find/Replace (rather than a
#    separate InsertBefore) keeps the edit scoped to that one run instead
#    of disturbing the whole paragraph.
$rngMps = $d.Content
$rngMps.Find.Execute("_mps_", $false, $false, $false, $false, $false, $true, 1, $false, " - _mps_", 2)

# 2) Drop the bookmark from its old spot (right after " County", before
#    ", " in the Location paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3) Re-create the "_GoBack" bookmark immediately before "_mps_", i.e. in
#    its new home right after the " - " we just inserted.
$rngTarget = $d.Content
$rngTarget.Find.Execute("_mps_")
$bmPoint = $rngTarget.Duplicate
$bmPoint.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmPoint)
